$wb = $excel.ActiveWorkbook

# Sheets 1,2,3 ("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")
# each had a row 2 for "QUIMICA I"/"1AM" that is being removed entirely,
# causing the former row 3 ("LECTURA.../1AV") to shift up to row 2 with
# updated totals.

$names = @("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")
foreach ($name in $names) {
    $ws = $wb.Worksheets.Item($name)
    # Delete the whole row 2 (QUIMICA I / 1AM), shifting row 3 up to row 2
    $ws.Rows.Item(2).Delete()
}

# Update the new row 2 values (previously row 3) on each sheet
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("C2").Value = 40
$ws1.Range("D2").Value = 12
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 28
$ws1.Range("G2").Value = 70
$ws1.Range("H2").Value = 8.300000000000001

$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("C2").Value = 40
$ws2.Range("D2").Value = 40
$ws2.Range("E2").Value = 28
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0

$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("C2").Value = 40
$ws3.Range("D2").Value = 12
$ws3.Range("E2").Value = 0
$ws3.Range("F2").Value = 28
$ws3.Range("G2").Value = 70
$ws3.Range("H2").Value = 8.300000000000001
